$d = $word.ActiveDocument

# The document contains three occurrences of a run-split pattern:
#   <id>  |  p130v_N  |  </id>
# each as a separate w:r. The edit merges each triple of runs into a
# single run (keeping the formatting of the leading "<id>" run) whose
# text is the concatenation "<id>p130v_N</id>".
#
# Find.Execute, when given matching Find/Replace text across a range
# that spans multiple runs, collapses the match into a single run using
# the formatting of the first run in the matched range - which is
# exactly the desired transformation here.

$ids = @("p130v_1", "p130v_2", "p130v_3")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
    Write-Host "Replaced $needle : $found"
}
